$d = $word.ActiveDocument

# The target paragraph block lives inside a legacy VML text box (w:pict/v:shape)
# which is not reachable through Document.Paragraphs / Document.Content as text
# (VML text-box content is outside the normal "story" ranges), so we operate on
# the canonical WordOpenXML for the whole document, replace the exact
# paragraph-run sequence for the "FOR w IN work" block, and write the whole XML
# back with InsertXML (the only supported way to mutate such content).

$xml = $d.Content.WordOpenXML

# --- locate the start paragraph ("+++ FOR w IN work+++") ---
$startTextMarker = "+++ FOR w IN work+++"
$startTextIdx = $xml.IndexOf($startTextMarker)
if ($startTextIdx -lt 0) {
    throw "Could not locate start marker text for work block"
}
$beforeStart = $xml.Substring(0, $startTextIdx)
$startIdx = $beforeStart.LastIndexOf("<w:p ")
if ($startIdx -lt 0) {
    throw "Could not locate start paragraph open tag for work block"
}

# --- locate the end paragraph ("+++END-FOR w+++") ---
$endTextMarker = "+++END-FOR w+++"
$endTextIdx = $xml.IndexOf($endTextMarker)
if ($endTextIdx -lt 0) {
    throw "Could not locate end marker text for work block"
}
$closeTag = "</w:p>"
$endCloseIdx = $xml.IndexOf($closeTag, $endTextIdx)
if ($endCloseIdx -lt 0) {
    throw "Could not locate end paragraph close tag for work block"
}
$endIdx = $endCloseIdx + $closeTag.Length

$newBlock = (
    '<w:p><w:r><w:t>+++HTML `</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>&lt;</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>meta</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>charset</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="UTF-8"&gt;</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>&lt;</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>body</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve"> ${</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>work }</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>&lt;/body&gt;</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>`+++</w:t></w:r></w:p>' +
    '<w:p/>'
)

$newXml = $xml.Substring(0, $startIdx) + $newBlock + $xml.Substring($endIdx)

$d.Content.InsertXML($newXml)
